# edit.ps1 -- add "rho_max" sheet, populate min/max purity data
# (computeRho.m min/max purity cases; Bloch-sphere Stokes-parameter plot data)

$wb = $excel.ActiveWorkbook

$headers = @("theta", "Jxx", "Jyy", "beta", "gamma", "trace_sq")

# ------------------------------------------------------------------
# rho_mat: header formatting refreshed (values unchanged)
# ------------------------------------------------------------------
$wsMat = $wb.Worksheets.Item("rho_mat")
$wsMat.Range("A1:F1").NumberFormat = "@"

# ------------------------------------------------------------------
# rho_min: refreshed min-purity values (computeRho.m min-purity case)
# ------------------------------------------------------------------
$wsMin = $wb.Worksheets.Item("rho_min")
$wsMin.Range("A1:F1").NumberFormat = "@"
for ($col = 1; $col -le 6; $col++) {
  $wsMin.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$rhoMinData = @(
    @(0, 0.11740910107661572, 0.88259089892338438, -0.0024460142171670498, 0.00097044975456231711, 1),
    @(5, 0.14291919417730423, 0.85708080582269575, -0.14579038817613021, -0.005551633035313322, 1),
    @(10, 0.21320731585731345, 0.78679268414268655, -0.26693398702951754, -0.0015297245111286801, 1),
    @(15, 0.3164164511280857, 0.6835835488719143, -0.35459288206766981, 0.0015701812468705552, 1),
    @(20, 0.44366616885138704, 0.5563338311486129, -0.38439555371994871, 0.0040370616529314327, 1),
    @(25, 0.58745842508170998, 0.41254157491829002, -0.38233740620011769, 0.003798250835471913, 1),
    @(30, 0.70102127602277386, 0.29897872397722614, -0.32519134635403596, 0.0045262707761032335, 1),
    @(35, 0.80698263373966461, 0.19301736626033536, -0.23921604792516965, 0.0021981141287036815, 1),
    @(40, 0.86583343756364206, 0.13416656243635783, -0.12731678288609377, -0.0073147125341848293, 0.99999999999999978),
    @(45, 0.88228289107559965, 0.1177171089244004, 0.010427128095225294, -0.005761075915754644, 1),
    @(50, 0.85842762164147413, 0.14157237835852585, 0.13468917328317698, -0.014573726512641237, 1),
    @(55, 0.78756159961208505, 0.21243840038791484, 0.25709376346270968, -0.012998312576132475, 0.99999999999999978),
    @(60, 0.68624585663153714, 0.31375414336846291, 0.33200348356056542, -0.020199949053652258, 1),
    @(65, 0.55958963020258734, 0.44041036979741255, 0.38120190525563014, -0.021216694356277502, 0.99999999999999978),
    @(70, 0.422221007417727, 0.57777899258227305, 0.3760068625814737, -0.025212904642143461, 1),
    @(75, 0.30361227369827798, 0.69638772630172208, 0.33096276896531135, -0.022541234480017162, 1),
    @(80, 0.20145785725907406, 0.79854214274092594, 0.24572177558532818, -0.02290011836432726, 1),
    @(85, 0.1383262958931758, 0.8616737041068242, 0.12522929546073336, -0.019354090984365371, 1),
    @(90, 0.11553731031223477, 0.88446268968776531, -0.0081182864134879181, -0.017924891845677256, 1)
  )
for ($i = 0; $i -lt $rhoMinData.Count; $i++) {
  $r = $rhoMinData[$i]
  for ($col = 1; $col -le 6; $col++) {
    $wsMin.Cells.Item($i + 2, $col).Value = $r[$col - 1]
  }
}

# ------------------------------------------------------------------
# rho_max: new sheet, max-purity values (computeRho.m max-purity case)
# ------------------------------------------------------------------
$wsMax = $wb.Worksheets.Add($null, $wsMin)
$wsMax.Name = "rho_max"
$wsMax.Range("A1:F1").NumberFormat = "@"
for ($col = 1; $col -le 6; $col++) {
  $wsMax.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$rhoMaxData = @(
    @(0, 0.000011827304300320885, 0.99998817269569962, -0.0031965689258958642, 0.0012682303756884966, 0.99999999999999978),
    @(5, 0.03714368425529993, 0.9628563157447001, -0.18897683897268816, -0.0071961538430256556, 1),
    @(10, 0.13400495918802599, 0.86599504081197409, -0.34065204894963341, -0.00195218224117275, 1),
    @(15, 0.27011899708250292, 0.72988100291749713, -0.44401673166256245, 0.0019661611403137577, 1),
    @(20, 0.42750248730676937, 0.57249751269323068, -0.49468891014204658, 0.0051954025220590323, 1),
    @(25, 0.61148836508285787, 0.38851163491714225, -0.4873878335614571, 0.0048418522906825126, 1),
    @(30, 0.76288685087049535, 0.23711314912950462, -0.42527104923791237, 0.0059192593642780744, 1),
    @(35, 0.89438837203517929, 0.10561162796482065, -0.30732692125478711, 0.002823972946633039, 1),
    @(40, 0.97213603276519445, 0.027863967234805654, -0.16431204642361522, -0.0094401960075261616, 1),
    @(45, 0.99975740398089963, 0.00024259601910047343, 0.013631356750453081, -0.0075314391802720013, 1),
    @(50, 0.967706124132382, 0.032293875867617998, 0.17575361773284789, -0.01901700853906527, 1),
    @(55, 0.87253683609488475, 0.12746316390511536, 0.33306567131816489, -0.016839349371463273, 1),
    @(60, 0.74428181299466611, 0.25571818700533389, 0.43545888403396932, -0.026494442703164589, 1),
    @(65, 0.57710589309753924, 0.42289410690246076, 0.4932555086395326, -0.027453302887710553, 1),
    @(70, 0.39893424488176932, 0.60106575511823068, 0.48858202240452442, -0.032761561467728989, 1),
    @(75, 0.24528382253392447, 0.75471617746607556, 0.4292608961973241, -0.029236129926446055, 1),
    @(80, 0.11462359404122846, 0.88537640595877154, 0.31719265451597567, -0.029560869464695157, 1),
    @(85, 0.028123835270552448, 0.9718761647294476, 0.16338682901958501, -0.025251308352872832, 1),
    @(90, 0.00065361875764962112, 0.99934638124235042, -0.010544162154606953, -0.023281140452324864, 1)
  )
for ($i = 0; $i -lt $rhoMaxData.Count; $i++) {
  $r = $rhoMaxData[$i]
  for ($col = 1; $col -le 6; $col++) {
    $wsMax.Cells.Item($i + 2, $col).Value = $r[$col - 1]
  }
}

